# Auto-generated script to apply numeric cell updates per commit diff
# "chore: update Sheets via scheduled runner" -- refreshed market-board price snapshots
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2578.682
$ws.Range("I15").Value = 2578.682
$ws.Range("K15").Value = 7736.045999999999
$ws.Range("M15").Value = -7567.045999999999
$ws.Range("H42").Value = 42.3
$ws.Range("I42").Value = 46
$ws.Range("J42").Value = 27.5
$ws.Range("K42").Value = 138
$ws.Range("L42").Value = 82.5
$ws.Range("M42").Value = 92
$ws.Range("N42").Value = -542.5
$ws.Range("H43").Value = 5136.2144
$ws.Range("J43").Value = 5336.6
$ws.Range("L43").Value = 5336.6
$ws.Range("N43").Value = -5474.6
$ws.Range("H86").Value = 13719.77
$ws.Range("I86").Value = 12435.7
$ws.Range("K86").Value = 12435.7
$ws.Range("M86").Value = -11312.7
$ws.Range("H89").Value = 13719.77
$ws.Range("I89").Value = 12435.7
$ws.Range("K89").Value = 62178.5
$ws.Range("M89").Value = -56562.5
$ws.Range("H112").Value = 2112.353
$ws.Range("I112").Value = 2212.7144
$ws.Range("K112").Value = 6638.1432
$ws.Range("M112").Value = -5530.1432
$ws.Range("H113").Value = 4847.1
$ws.Range("I113").Value = 3998.8
$ws.Range("J113").Value = 5695.4
$ws.Range("K113").Value = 3998.8
$ws.Range("L113").Value = 5695.4
$ws.Range("M113").Value = -744.8000000000002
$ws.Range("N113").Value = -12203.4
$ws.Range("H132").Value = 1099.6285
$ws.Range("I132").Value = 866.3333
$ws.Range("K132").Value = 2598.9999
$ws.Range("M132").Value = -68.9998999999998
$ws.Range("H133").Value = 99000
$ws.Range("J133").Value = 99000
$ws.Range("L133").Value = 99000
$ws.Range("N133").Value = -109120
$ws.Range("H134").Value = 150000
$ws.Range("J134").Value = 150000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -160140
$ws.Range("H137").Value = 14290141
$ws.Range("I137").Value = 16130688
$ws.Range("J137").Value = 25902.25
$ws.Range("K137").Value = 48392064
$ws.Range("L137").Value = 77706.75
$ws.Range("M137").Value = -48389514
$ws.Range("N137").Value = -82806.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2299.3125
$ws.Range("I45").Value = 1837.6154
$ws.Range("J45").Value = 4300
$ws.Range("K45").Value = 1837.6154
$ws.Range("L45").Value = 4300
$ws.Range("M45").Value = -1460.6154
$ws.Range("N45").Value = -5054
$ws.Range("H74").Value = 1510403.2
$ws.Range("I74").Value = 2228699.5
$ws.Range("J74").Value = 13952.583
$ws.Range("K74").Value = 2228699.5
$ws.Range("L74").Value = 13952.583
$ws.Range("M74").Value = -2227825.5
$ws.Range("N74").Value = -15700.583
$ws.Range("H77").Value = 1510403.2
$ws.Range("I77").Value = 2228699.5
$ws.Range("J77").Value = 13952.583
$ws.Range("K77").Value = 11143497.5
$ws.Range("L77").Value = 69762.91500000001
$ws.Range("M77").Value = -11139129.5
$ws.Range("N77").Value = -78498.91500000001
$ws.Range("H122").Value = 2208.476
$ws.Range("I122").Value = 1675.2307
$ws.Range("K122").Value = 5025.6921
$ws.Range("M122").Value = -2575.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 32552.5
$ws.Range("I75").Value = 25070
$ws.Range("K75").Value = 25070
$ws.Range("M75").Value = -24134
$ws.Range("H78").Value = 32552.5
$ws.Range("I78").Value = 25070
$ws.Range("K78").Value = 75210
$ws.Range("M78").Value = -70530
$ws.Range("H82").Value = 16929.143
$ws.Range("I82").Value = 4351.2
$ws.Range("J82").Value = 48374
$ws.Range("K82").Value = 4351.2
$ws.Range("L82").Value = 48374
$ws.Range("M82").Value = -3968.2
$ws.Range("N82").Value = -49140
$ws.Range("H85").Value = 16929.143
$ws.Range("I85").Value = 4351.2
$ws.Range("J85").Value = 48374
$ws.Range("K85").Value = 4351.2
$ws.Range("L85").Value = 48374
$ws.Range("M85").Value = -3025.2
$ws.Range("N85").Value = -51026
$ws.Range("H107").Value = 1582.8055
$ws.Range("J107").Value = 1919.1875
$ws.Range("L107").Value = 1919.1875
$ws.Range("N107").Value = -5759.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3217.3076
$ws.Range("I16").Value = 3872.6
$ws.Range("J16").Value = 1033
$ws.Range("K16").Value = 3872.6
$ws.Range("L16").Value = 1033
$ws.Range("M16").Value = -3585.6
$ws.Range("N16").Value = -1607
$ws.Range("H68").Value = 38925
$ws.Range("J68").Value = 38925
$ws.Range("L68").Value = 38925
$ws.Range("N68").Value = -40423
$ws.Range("H71").Value = 38925
$ws.Range("J71").Value = 38925
$ws.Range("L71").Value = 116775
$ws.Range("N71").Value = -124263
$ws.Range("H75").Value = 45000
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null
$ws.Range("H78").Value = 45000
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null
$ws.Range("H99").Value = 36983.168
$ws.Range("I99").Value = 71299.664
$ws.Range("J99").Value = 2666.6667
$ws.Range("K99").Value = 71299.664
$ws.Range("L99").Value = 2666.6667
$ws.Range("M99").Value = -69801.664
$ws.Range("N99").Value = -5662.6667
$ws.Range("H113").Value = 3217.3076
$ws.Range("I113").Value = 3872.6
$ws.Range("J113").Value = 1033
$ws.Range("K113").Value = 3872.6
$ws.Range("L113").Value = 1033
$ws.Range("M113").Value = -1702.6
$ws.Range("N113").Value = -5373
$ws.Range("H126").Value = 36983.168
$ws.Range("I126").Value = 71299.664
$ws.Range("J126").Value = 2666.6667
$ws.Range("K126").Value = 213898.992
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("M126").Value = -211428.992
$ws.Range("N126").Value = -12940.0001
$ws.Range("H132").Value = 9077.8125
$ws.Range("I132").Value = 9393
$ws.Range("K132").Value = 28179
$ws.Range("M132").Value = -25649
$ws.Range("H134").Value = 3009.3684
$ws.Range("I134").Value = 2621
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 7863
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -5328
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2388.4707
$ws.Range("J34").Value = 2533.1428
$ws.Range("L34").Value = 7599.428400000001
$ws.Range("N34").Value = -7767.428400000001
$ws.Range("H97").Value = 891.6
$ws.Range("J97").Value = 891.6
$ws.Range("L97").Value = 2674.8
$ws.Range("N97").Value = -3666.8
$ws.Range("H98").Value = 494.25
$ws.Range("J98").Value = 374.75
$ws.Range("L98").Value = 1124.25
$ws.Range("N98").Value = -4120.25
$ws.Range("H107").Value = 3375.1177
$ws.Range("J107").Value = 4167.25
$ws.Range("L107").Value = 12501.75
$ws.Range("N107").Value = -16341.75
$ws.Range("H123").Value = 11666.5
$ws.Range("I123").Value = 8333
$ws.Range("K123").Value = 24999
$ws.Range("M123").Value = -22549

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null
$ws.Range("H102").Value = 10350
$ws.Range("I102").Value = 10186.454
$ws.Range("K102").Value = 10186.454
$ws.Range("M102").Value = -8564.454
$ws.Range("H132").Value = 12916.556
$ws.Range("I132").Value = 12383.444
$ws.Range("J132").Value = 13449.667
$ws.Range("K132").Value = 37150.33199999999
$ws.Range("L132").Value = 40349.001
$ws.Range("M132").Value = -34620.33199999999
$ws.Range("N132").Value = -45409.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3532.5
$ws.Range("I22").Value = 2874.25
$ws.Range("J22").Value = 3908.6428
$ws.Range("K22").Value = 2874.25
$ws.Range("L22").Value = 3908.6428
$ws.Range("M22").Value = -2579.25
$ws.Range("N22").Value = -4498.6428
$ws.Range("H27").Value = 3532.5
$ws.Range("I27").Value = 2874.25
$ws.Range("J27").Value = 3908.6428
$ws.Range("K27").Value = 2874.25
$ws.Range("L27").Value = 3908.6428
$ws.Range("M27").Value = -2767.25
$ws.Range("N27").Value = -4122.6428
$ws.Range("H46").Value = 3955.75
$ws.Range("J46").Value = 5143.4116
$ws.Range("L46").Value = 5143.4116
$ws.Range("N46").Value = -5519.4116
$ws.Range("H55").Value = 1368.5883
$ws.Range("J55").Value = 1272.4166
$ws.Range("L55").Value = 1272.4166
$ws.Range("N55").Value = -1618.4166
$ws.Range("H122").Value = 5580.1816
$ws.Range("I122").Value = 3277.4
$ws.Range("J122").Value = 7499.1665
$ws.Range("K122").Value = 9832.200000000001
$ws.Range("L122").Value = 22497.4995
$ws.Range("M122").Value = -7382.200000000001
$ws.Range("N122").Value = -27397.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47401.25
$ws.Range("I70").Value = 49833.332
$ws.Range("J70").Value = 40105
$ws.Range("K70").Value = 49833.332
$ws.Range("L70").Value = 40105
$ws.Range("M70").Value = -49518.332
$ws.Range("N70").Value = -40735
$ws.Range("H73").Value = 47401.25
$ws.Range("I73").Value = 49833.332
$ws.Range("J73").Value = 40105
$ws.Range("K73").Value = 49833.332
$ws.Range("L73").Value = 40105
$ws.Range("M73").Value = -48741.332
$ws.Range("N73").Value = -42289
$ws.Range("H88").Value = 14999.5
$ws.Range("I88").Value = 14999.5
$ws.Range("K88").Value = 14999.5
$ws.Range("M88").Value = -14593.5
$ws.Range("H91").Value = 14999.5
$ws.Range("I91").Value = 14999.5
$ws.Range("K91").Value = 14999.5
$ws.Range("M91").Value = -13595.5
